$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Marrks"
$ws.Range("A2").Value = "KPK"
$ws.Range("B2").Value = 96
$ws.Range("A3").Value = "Shruti"
$ws.Range("B3").Value = 98

$ws.Range("C5").Select()
